$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Email faktura/ordre/kontoudtog/provisionsseddel"
$ws.Range("A11").EntireRow.Delete()

# Update complexity ("Kompleksitet") values for a few functions
$ws.Range("B3").Value = "Kompleks"   # Opret, rediger, fjern vare: Simpel -> Kompleks
$ws.Range("B6").Value = "Kompleks"   # Opret, rediger, fjern ordre: Medium -> Kompleks
$ws.Range("B8").Value = "Medium"     # Opgrader til faktura: Simpel -> Medium

# Update the active selection
[void]$ws.Range("B4").Select()
